$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "About" sheet: insert a new "medium" header row (row 2, pushing the
#    "small" header down to row 3) and a new "medium" footnote row (pushing
#    the "small" footnote down), documenting the new BpTPEU-medium tab added
#    for the clean industrial heat PTC policy. Doing this first establishes
#    the shared-string insertion order to match the target workbook.
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

$about.Rows("2:2").Insert()
$about.Range("A2").Value = "BpTPEU BTU per Medium Primary Energy Unit"

$about.Rows("13:13").Insert()
$about.Range("A13").Value = "The medium primary energy output unit (used in the clean heat PTC policy) is: million BTU"

# ---------------------------------------------------------------------------
# 2. Add the new "BpTPEU-medium" worksheet by copying the existing
#    "BpTPEU-small" sheet (so it inherits the same tab color, column widths
#    and cell styles) and positioning the copy right after "BpTPEU-large".
# ---------------------------------------------------------------------------
$largeSheet = $wb.Worksheets.Item("BpTPEU-large")
$smallSheetOrig = $wb.Worksheets.Item("BpTPEU-small")

$largeIndex = $largeSheet.Index
$smallSheetOrig.Copy($null, $largeSheet)

# Sheet references in this host are resolved positionally, so re-fetch by
# position now that the sheet collection has shifted - the copy lands
# immediately after "BpTPEU-large".
$medium = $wb.Worksheets.Item($largeIndex + 1)
$medium.Name = "BpTPEU-medium"

$medium.Range("B1").Value = "medium primary energy output unit"
$medium.Range("B2").Formula = "=10^6"
$medium.Range("B2").Style = "Normal"

# ---------------------------------------------------------------------------
# 3. The original small sheet's B2 cell loses its now-unused custom style,
#    becoming an ordinary, unstyled numeric formula cell (re-fetch by name,
#    since the Copy above shifted what the old reference resolves to).
# ---------------------------------------------------------------------------
$smallSheet = $wb.Worksheets.Item("BpTPEU-small")
$smallSheet.Range("B2").Style = "Normal"

# ---------------------------------------------------------------------------
# 4. Restore "About" as the selected/active sheet.
# ---------------------------------------------------------------------------
$about.Activate()

$wb.Application.Calculate()
